$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously empty "Aufwand" (effort) cells for the second table (rows 14-20) ---
$ws.Range("C14").Value = "0.5 Stunden"
$ws.Range("C15").Value = "2 Stunden"
$ws.Range("C16").Value = "3 Stunden"
$ws.Range("C17").Value = "1.5 Stunden"
$ws.Range("C18").Value = "8 Stunden"
$ws.Range("C19").Value = "0.5 Stunden"

# Row 20 ("Release") gets a responsible person and an effort placeholder
$ws.Range("B20").Value = "Geschäftsführer"
$ws.Range("C20").Value = " --- "

# Row 21 ("Übergabe des Projekts") gets an effort value and its due date corrected
$ws.Range("B21").Copy($ws.Range("C21"))
$ws.Range("C21").Value = "0.25 Stunden"
$ws.Range("D21").Value = 43609

# --- New row 22: documentation/report update entry ---
# Carry over the same row formatting used by the rest of the table (row 21)
# before filling in the new row's content.
$ws.Range("A21:C21").Copy($ws.Range("A22:C22"))
$ws.Range("E21").Copy($ws.Range("E22"))

$ws.Range("A22").Value = "Update"
$ws.Range("B22").Value = "Architekt / Programmierer"
$ws.Range("C22").Value = "2.5 Stunden"
$ws.Range("D22").HorizontalAlignment = -4131
$ws.Range("D22").NumberFormat = "mm-dd-yy"
$ws.Range("D22").Value = 43637
$ws.Range("E22").Value = "Update und Dokument-Aktualisierungen"

# Move the active selection past the newly entered data, like a user would
# after finishing data entry in the last row.
$ws.Range("E23").Select()
